$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace static German header/labels with template placeholders
# (KIBON-120: Mitarbeiterinnen statistic translated)
$ws.Range("A1").Value = "{mitarbeiterinnenTitle}"

$ws.Range("A3").Value = "{vonTitle}"
$ws.Range("B3").Value = "{auswertungVon}"

$ws.Range("A4").Value = "{bisTitle}"
$ws.Range("B4").Value = "{auswertungBis}"

$ws.Range("A6").Value = "{nachnameTitle}"
$ws.Range("B6").Value = "{vornameTitle}"
$ws.Range("C6").Value = "{anzahlVerGesucheTitle}"
$ws.Range("D6").Value = "{verfuegungAusgestelltTitle}"

$ws.Range("A7").Value = "{name}"
$ws.Range("B7").Value = "{vorname}"
$ws.Range("C7").Value = "{verantwortlicheGesuche}"
$ws.Range("D7").Value = "{verfuegungenAusgestellt}"
$ws.Range("E7").Value = "{repeatMitarbeiterinnenRow}"

# Update the active selection recorded in the sheet view
$ws.Range("A6").Select()
